# Auto-generated edit script: updates crypto price/volume table
# Source: diff of cryptos.xlsx (GitHub Actions scheduled data refresh)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    # Force the cell to keep a literal text value (e.g. "1.00", "10.50")
    # instead of Excel auto-coercing it to a number and dropping trailing zeros.
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# Row 2
$ws.Cells.Item(2, 4).Value = "69.863.65"
$ws.Cells.Item(2, 5).Value = "  +0.56%  "

# Row 3
$ws.Cells.Item(3, 4).Value = "3.690.76"
$ws.Cells.Item(3, 5).Value = "  +0.00%  "

# Row 4
$ws.Cells.Item(4, 5).Value = "  -0.01%  "

# Row 5
Set-TextValue $ws.Cells.Item(5, 4) "648.81"
$ws.Cells.Item(5, 5).Value = "  -4.71%  "

# Row 6
Set-TextValue $ws.Cells.Item(6, 4) "161.01"
$ws.Cells.Item(6, 5).Value = "  -0.27%  "

# Row 7
$ws.Cells.Item(7, 5).Value = "  +0.02%  "

# Row 8
$ws.Cells.Item(8, 5).Value = "  +1.55%  "

# Row 9
$ws.Cells.Item(9, 5).Value = "  -0.38%  "

# Row 10
Set-TextValue $ws.Cells.Item(10, 4) "7.20"
$ws.Cells.Item(10, 5).Value = "  +0.75%  "

# Row 11
Set-TextValue $ws.Cells.Item(11, 4) "0.445"
$ws.Cells.Item(11, 5).Value = "  +1.12%  "

# Row 12
$ws.Cells.Item(12, 5).Value = "  -0.56%  "

# Row 13
$ws.Cells.Item(13, 4).Value = "4.314.14"
$ws.Cells.Item(13, 5).Value = "  -0.05%  "

# Row 14
Set-TextValue $ws.Cells.Item(14, 4) "32.80"
$ws.Cells.Item(14, 5).Value = "  +1.07%  "

# Row 15
$ws.Cells.Item(15, 4).Value = "3.689.15"
$ws.Cells.Item(15, 5).Value = "  -0.28%  "

# Row 16
$ws.Cells.Item(16, 4).Value = "69.846.36"
$ws.Cells.Item(16, 5).Value = "  +0.52%  "

# Row 17
$ws.Cells.Item(17, 5).Value = "  +0.31%  "

# Row 18
Set-TextValue $ws.Cells.Item(18, 4) "16.11"
$ws.Cells.Item(18, 5).Value = "  +0.53%  "

# Row 19
Set-TextValue $ws.Cells.Item(19, 4) "6.52"
$ws.Cells.Item(19, 5).Value = "  +0.70%  "

# Row 20
Set-TextValue $ws.Cells.Item(20, 4) "10.50"
$ws.Cells.Item(20, 5).Value = "  +7.27%  "

# Row 21
Set-TextValue $ws.Cells.Item(21, 4) "471.31"
$ws.Cells.Item(21, 5).Value = "  -0.01%  "

# Row 22
$ws.Cells.Item(22, 5).Value = "  +0.13%  "

# Row 23
Set-TextValue $ws.Cells.Item(23, 4) "79.92"
$ws.Cells.Item(23, 5).Value = "  -0.70%  "

# Row 24
$ws.Cells.Item(24, 4).Value = "3.837.29"
$ws.Cells.Item(24, 5).Value = "  -0.04%  "

# Row 25
$ws.Cells.Item(25, 5).Value = "  -0.01%  "

# Row 26
$ws.Cells.Item(26, 5).Value = "  -0.44%  "

# Row 27
Set-TextValue $ws.Cells.Item(27, 4) "10.94"
$ws.Cells.Item(27, 5).Value = "  +0.41%  "

# Row 28
Set-TextValue $ws.Cells.Item(28, 4) "9.16"
$ws.Cells.Item(28, 5).Value = "  +0.44%  "

# Row 29
$ws.Cells.Item(29, 5).Value = "  -2.00%  "

# Row 30
Set-TextValue $ws.Cells.Item(30, 4) "1.72"
$ws.Cells.Item(30, 5).Value = "  -1.14%  "

# Row 31
Set-TextValue $ws.Cells.Item(31, 4) "2.02"
$ws.Cells.Item(31, 5).Value = "  +0.14%  "

# Row 32
Set-TextValue $ws.Cells.Item(32, 4) "6.56"
$ws.Cells.Item(32, 5).Value = "  -0.22%  "

# Row 33
Set-TextValue $ws.Cells.Item(33, 4) "1.00"
$ws.Cells.Item(33, 5).Value = "  -0.20%  "

# Row 34
Set-TextValue $ws.Cells.Item(34, 4) "26.86"
$ws.Cells.Item(34, 5).Value = "  -0.49%  "

# Row 35
$ws.Cells.Item(35, 4).Value = "3.688.71"
$ws.Cells.Item(35, 5).Value = "  +0.16%  "

# Row 36
Set-TextValue $ws.Cells.Item(36, 4) "0.163"
$ws.Cells.Item(36, 5).Value = "  +0.27%  "

# Row 37
Set-TextValue $ws.Cells.Item(37, 4) "8.47"
$ws.Cells.Item(37, 5).Value = "  +0.20%  "

# Row 39
Set-TextValue $ws.Cells.Item(39, 4) "5.91"
$ws.Cells.Item(39, 5).Value = "  -4.84%  "

# Row 40
$ws.Cells.Item(40, 2).Value = "Stacks"
$ws.Cells.Item(40, 3).Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
Set-TextValue $ws.Cells.Item(40, 4) "2.26"
$ws.Cells.Item(40, 5).Value = "  -0.96%  "

# Row 41
$ws.Cells.Item(41, 2).Value = "Monero"
$ws.Cells.Item(41, 3).Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextValue $ws.Cells.Item(41, 4) "179.30"
$ws.Cells.Item(41, 5).Value = "  +6.37%  "

# Row 42
$ws.Cells.Item(42, 5).Value = "  -0.02%  "

# Row 43
$ws.Cells.Item(43, 5).Value = "  +0.41%  "

# Row 44
Set-TextValue $ws.Cells.Item(44, 4) "0.929"
$ws.Cells.Item(44, 5).Value = "  -1.41%  "

# Row 45
Set-TextValue $ws.Cells.Item(45, 4) "47.16"
$ws.Cells.Item(45, 5).Value = "  +0.92%  "

# Row 46
Set-TextValue $ws.Cells.Item(46, 4) "29.12"
$ws.Cells.Item(46, 5).Value = "  +5.56%  "

# Row 47
$ws.Cells.Item(47, 5).Value = "  -0.88%  "

# Row 48
$ws.Cells.Item(48, 5).Value = "  -0.65%  "

# Row 49
$ws.Cells.Item(49, 2).Value = "FLOKI"
$ws.Cells.Item(49, 3).Value = "https://coinranking.com/coin/fmHk13Rqw+floki-floki"
Set-TextValue $ws.Cells.Item(49, 4) "0.000268"
$ws.Cells.Item(49, 5).Value = "  -4.07%  "

# Row 50
$ws.Cells.Item(50, 2).Value = "Cosmos"
$ws.Cells.Item(50, 3).Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
Set-TextValue $ws.Cells.Item(50, 4) "7.85"
$ws.Cells.Item(50, 5).Value = "  -0.60%  "

# Row 51
Set-TextValue $ws.Cells.Item(51, 4) "1.25"
$ws.Cells.Item(51, 5).Value = "  -3.72%  "
